# Azerbaijan Premier League workbook update (19-06-2024 21:51)
# Rotates betting-data rows 172-174: row 174 moves to 172, 172 moves to 173, 173 moves to 174
# (column A, the sequential index, stays fixed per row)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 172 (content copied from current row 174)
$ws.Range("B172").Value = 7158118
$ws.Range("C172").Value = "Azerbaijan Premier League"
$ws.Range("D172").Value = 45430.41666666666
$ws.Range("E172").Value = "Sabail FC"
$ws.Range("F172").Value = "Neftchi Baku"
$ws.Range("G172").Value = 0
$ws.Range("H172").Value = 3
$ws.Range("I172").Value = 0
$ws.Range("J172").Value = 0
$ws.Range("K172").Value = "A"
$ws.Range("L172").Value = 4.5
$ws.Range("M172").Value = 3.7
$ws.Range("N172").Value = 1.6
$ws.Range("O172").Value = 4.2
$ws.Range("P172").Value = 3.5
$ws.Range("Q172").Value = 1.7
$ws.Range("R172").Value = 0.75
$ws.Range("S172").Value = 1.9
$ws.Range("T172").Value = 1.9
$ws.Range("U172").Value = 2.75
$ws.Range("V172").Value = 1.825
$ws.Range("W172").Value = 1.975
$ws.Range("X172").Value = -1
$ws.Range("Y172").Value = -1
$ws.Range("Z172").Value = 0.7
$ws.Range("AA172").Value = -1
$ws.Range("AB172").Value = 0.8999999999999999
$ws.Range("AC172").Value = 0.4125
$ws.Range("AD172").Value = -0.5

# Row 173 (content copied from current row 172)
$ws.Range("B173").Value = 7157321
$ws.Range("C173").Value = "Azerbaijan Premier League"
$ws.Range("D173").Value = 45430.41666666666
$ws.Range("E173").Value = "FK Kapaz"
$ws.Range("F173").Value = "Zira IK"
$ws.Range("G173").Value = 0
$ws.Range("H173").Value = 1
$ws.Range("I173").Value = 0
$ws.Range("J173").Value = 0
$ws.Range("K173").Value = "A"
$ws.Range("L173").Value = 4.5
$ws.Range("M173").Value = 3.4
$ws.Range("N173").Value = 1.666
$ws.Range("O173").Value = 5
$ws.Range("P173").Value = 3.6
$ws.Range("Q173").Value = 1.533
$ws.Range("R173").Value = 1
$ws.Range("S173").Value = 1.825
$ws.Range("T173").Value = 1.975
$ws.Range("U173").Value = 2.5
$ws.Range("V173").Value = 1.875
$ws.Range("W173").Value = 1.925
$ws.Range("X173").Value = -1
$ws.Range("Y173").Value = -1
$ws.Range("Z173").Value = 0.5329999999999999
$ws.Range("AA173").Value = 0
$ws.Range("AB173").Value = 0
$ws.Range("AC173").Value = -1
$ws.Range("AD173").Value = 0.925

# Row 174 (content copied from current row 173)
$ws.Range("B174").Value = 7153759
$ws.Range("C174").Value = "Azerbaijan Premier League"
$ws.Range("D174").Value = 45430.41666666666
$ws.Range("E174").Value = "Sabah"
$ws.Range("F174").Value = "FK Sumqayit"
$ws.Range("G174").Value = 2
$ws.Range("H174").Value = 0
$ws.Range("I174").Value = 1
$ws.Range("J174").Value = 0
$ws.Range("K174").Value = "H"
$ws.Range("L174").Value = 1.833
$ws.Range("M174").Value = 3.3
$ws.Range("N174").Value = 3.8
$ws.Range("O174").Value = 2.15
$ws.Range("P174").Value = 3.2
$ws.Range("Q174").Value = 3.1
$ws.Range("R174").Value = -0.25
$ws.Range("S174").Value = 1.9
$ws.Range("T174").Value = 1.9
$ws.Range("U174").Value = 2.25
$ws.Range("V174").Value = 1.8
$ws.Range("W174").Value = 2
$ws.Range("X174").Value = 1.15
$ws.Range("Y174").Value = -1
$ws.Range("Z174").Value = -1
$ws.Range("AA174").Value = 0.8999999999999999
$ws.Range("AB174").Value = -1
$ws.Range("AC174").Value = -0.5
$ws.Range("AD174").Value = 0.5

